$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Ja Morant",           "PG",       "Memphis Grizzlies"),
    @("Bennedict Mathurin",  "SG,SF",    "Indiana Pacers"),
    @("Caris LeVert",        "SG,SF",    "Cleveland Cavaliers"),
    @("Santi Aldama",        "PF,C",     "Memphis Grizzlies"),
    @("Nikola Vucevic",      "PF,C",     "Chicago Bulls"),
    @("Evan Mobley",         "PF,C",     "Cleveland Cavaliers"),
    @("Nick Richards",       "C",        "Charlotte Hornets"),
    @("Luka Doncic",         "PG,SG",    "Dallas Mavericks"),
    @("Brook Lopez",         "C",        "Milwaukee Bucks"),
    @("DeMar DeRozan",       "SF,PF",    "Sacramento Kings"),
    @("Gradey Dick",         "SG,SF",    "Toronto Raptors"),
    @("De'Aaron Fox",        "PG",       "Sacramento Kings"),
    @("Tyler Herro",         "PG,SG",    "Miami Heat"),
    @("Amen Thompson",       "SG,SF",    "Houston Rockets"),
    @("Mikal Bridges",       "SG,SF,PF", "New York Knicks"),
    @("Scottie Barnes",      "SG,SF,PF", "Toronto Raptors"),
    @("Miles Bridges",       "SF,PF",    "Charlotte Hornets"),
    @("Josh Giddey",         "PG,SG,SF", "Chicago Bulls")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
